$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.127881588408715, 0.3127903958511391, 0.8054896365839992, 8.660232485948974, 9.906394106792828)
    3 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182)
    4 = @(0.04763786555579896, 0.3127903958511391, 3.900430680208489, 0.496779210170732, 4.757638151786159)
    5 = @(0.127881588408715, 1.667794583268128, 0.8054896365839992, 8.660232485948974, 11.26139829420982)
    6 = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144)
    7 = @(0.3048080303191223, 1.667794583268128, 0.1575252929769615, 8.660232485948974, 10.79036039251319)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 7).Value = $vals[4]
}
